# Auto-generated edit script: updates cryptos.xlsx "Price" (D) and "Volume(1h)" (E) columns
# with the latest scraped values (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A number of the new Price values are purely numeric-looking strings (e.g. "1.002").
# Pre-format those cells as Text so Excel keeps them as literal strings instead of
# silently converting them to numbers (the source data models Price as text, same as
# the untouched cells that already contain non-numeric-looking strings like "29.867.76").
# NumberFormat is applied per contiguous block (not as one multi-area Range) since the
# host only reliably keeps the Text format on the first area of a multi-area Range.
$ws.Range("D4:D6").NumberFormat = "@"
$ws.Range("D8:D12").NumberFormat = "@"
$ws.Range("D14:D16").NumberFormat = "@"
$ws.Range("D19:D20").NumberFormat = "@"
$ws.Range("D23:D28").NumberFormat = "@"
$ws.Range("D30:D41").NumberFormat = "@"
$ws.Range("D43:D45").NumberFormat = "@"
$ws.Range("D47:D50").NumberFormat = "@"

$ws.Range("D2").Value = "29.867.76"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.886.49"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "0.7498"
$ws.Range("E5").Value = "  -5.14%  "
$ws.Range("D6").Value = "242.09"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").Value = "0.3113"
$ws.Range("E8").Value = "  -1.58%  "
$ws.Range("D9").Value = "25.39"
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("D10").Value = "0.07123"
$ws.Range("E10").Value = "  -2.90%  "
$ws.Range("D11").Value = "0.08536"
$ws.Range("D12").Value = "0.7593"
$ws.Range("E12").Value = "  -2.46%  "
$ws.Range("D13").Value = "1.909.04"
$ws.Range("E13").Value = "  +4.52%  "
$ws.Range("D14").Value = "5.357"
$ws.Range("E14").Value = "  -3.13%  "
$ws.Range("D15").Value = "93.36"
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("D16").Value = "6.135"
$ws.Range("E16").Value = "  -2.09%  "
$ws.Range("D17").Value = "29.895.58"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("E18").Value = "  -2.17%  "
$ws.Range("D19").Value = "243.34"
$ws.Range("E19").Value = "  -1.72%  "
$ws.Range("D20").Value = "0.000007794"
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("D21").Value = "2.152.67"
$ws.Range("E21").Value = "  +4.07%  "
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").Value = "7.981"
$ws.Range("E23").Value = "  -2.33%  "
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("D25").Value = "0.1602"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("D26").Value = "9.357"
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("D27").Value = "162.60"
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("D28").Value = "18.72"
$ws.Range("E28").Value = "  -0.76%  "
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("D30").Value = "1.522"
$ws.Range("E30").Value = "  +5.09%  "
$ws.Range("D31").Value = "1.534"
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("D32").Value = "4.469"
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("D33").Value = "4.094"
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("D34").Value = "0.05405"
$ws.Range("E34").Value = "  -4.08%  "
$ws.Range("D35").Value = "1.237"
$ws.Range("E35").Value = "  -1.27%  "
$ws.Range("D36").Value = "0.7429"
$ws.Range("E36").Value = "  -1.88%  "
$ws.Range("D37").Value = "0.9997"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "2.712"
$ws.Range("E38").Value = "  +1.69%  "
$ws.Range("D39").Value = "0.01938"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").Value = "2.779"
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("D41").Value = "0.4448"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("D42").Value = "1.103.71"
$ws.Range("E42").Value = "  -3.93%  "
$ws.Range("D43").Value = "6.067"
$ws.Range("E43").Value = "  +1.34%  "
$ws.Range("D44").Value = "72.43"
$ws.Range("E44").Value = "  -2.18%  "
$ws.Range("D45").Value = "0.8581"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D47").Value = "102.48"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("D48").Value = "1.864"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("D49").Value = "7.640"
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("D50").Value = "3.054"
$ws.Range("E50").Value = "  -3.88%  "
$ws.Range("D51").Value = "2.044.28"
$ws.Range("E51").Value = "  +2.60%  "
